# This script applies a row-content rearrangement to the "Artfynd" sheet:
#  - Rows 3 and 4 swap their full content.
#  - Rows 25, 26 and 27 rotate (new25=old27, new26=old25, new27=old26).
# Only cells whose value actually changes are written, so cells that are
# identical between source and destination (including cells that are blank
# both before and after) are left completely untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# Row 3 <-> Row 4 full swap
# ---------------------------------------------------------------------

# --- capture current (pre-edit) values needed from row 3 ---
$r3_A  = $ws.Range("A3").Value2
$r3_B  = $ws.Range("B3").Value2
$r3_E  = $ws.Range("E3").Value2
$r3_F  = $ws.Range("F3").Value2
$r3_G  = $ws.Range("G3").Value2
$r3_H  = $ws.Range("H3").Value2
$r3_Q  = $ws.Range("Q3").Value2
$r3_R  = $ws.Range("R3").Value2
$r3_AC = $ws.Range("AC3").Value2

# --- capture current (pre-edit) values needed from row 4 ---
$r4_A  = $ws.Range("A4").Value2
$r4_B  = $ws.Range("B4").Value2
$r4_E  = $ws.Range("E4").Value2
$r4_F  = $ws.Range("F4").Value2
$r4_G  = $ws.Range("G4").Value2
$r4_K  = $ws.Range("K4").Value2
$r4_Q  = $ws.Range("Q4").Value2
$r4_R  = $ws.Range("R4").Value2
$r4_AC = $ws.Range("AC4").Value2
$r4_AH = $ws.Range("AH4").Value2
$r4_AJ = $ws.Range("AJ4").Value2
$r4_AK = $ws.Range("AK4").Value2
$r4_AM = $ws.Range("AM4").Value2
$r4_AO = $ws.Range("AO4").Value2

# --- write row 3 with row 4's former values (only changed cells) ---
$ws.Range("A3").Value = $r4_A
$ws.Range("B3").Value = $r4_B
$ws.Range("E3").Value = $r4_E
$ws.Range("F3").Value = $r4_F
$ws.Range("G3").Value = $r4_G
$ws.Range("H3").ClearContents()
$ws.Range("K3").Value = $r4_K
$ws.Range("Q3").Value = $r4_Q
$ws.Range("R3").Value = $r4_R
$ws.Range("AC3").Value = $r4_AC
$ws.Range("AH3").Value = $r4_AH
$ws.Range("AJ3").Value = $r4_AJ
$ws.Range("AK3").Value = $r4_AK
$ws.Range("AM3").Value = $r4_AM
$ws.Range("AO3").Value = $r4_AO

# --- write row 4 with row 3's former values (only changed cells) ---
$ws.Range("A4").Value = $r3_A
$ws.Range("B4").Value = $r3_B
$ws.Range("E4").Value = $r3_E
$ws.Range("F4").Value = $r3_F
$ws.Range("G4").Value = $r3_G
$ws.Range("H4").Value = $r3_H
$ws.Range("K4").ClearContents()
$ws.Range("Q4").Value = $r3_Q
$ws.Range("R4").Value = $r3_R
$ws.Range("AC4").Value = $r3_AC
$ws.Range("AH4").ClearContents()
$ws.Range("AJ4").ClearContents()
$ws.Range("AK4").ClearContents()
$ws.Range("AM4").ClearContents()
$ws.Range("AO4").ClearContents()

# ---------------------------------------------------------------------
# Rows 25, 26, 27 rotation: new25=old27, new26=old25, new27=old26
# ---------------------------------------------------------------------

# --- capture current (pre-edit) values ---
$r25_A = $ws.Range("A25").Value2
$r25_B = $ws.Range("B25").Value2
$r25_E = $ws.Range("E25").Value2
$r25_F = $ws.Range("F25").Value2
$r25_G = $ws.Range("G25").Value2
$r25_H = $ws.Range("H25").Value2
$r25_Q = $ws.Range("Q25").Value2
$r25_R = $ws.Range("R25").Value2
$r25_S = $ws.Range("S25").Value2
$r25_Z = $ws.Range("Z25").Value2
$r25_AB = $ws.Range("AB25").Value2

$r26_A = $ws.Range("A26").Value2
$r26_B = $ws.Range("B26").Value2
$r26_E = $ws.Range("E26").Value2
$r26_F = $ws.Range("F26").Value2
$r26_G = $ws.Range("G26").Value2
$r26_H = $ws.Range("H26").Value2
$r26_Q = $ws.Range("Q26").Value2
$r26_R = $ws.Range("R26").Value2
$r26_S = $ws.Range("S26").Value2
$r26_Z = $ws.Range("Z26").Value2
$r26_AB = $ws.Range("AB26").Value2

$r27_A = $ws.Range("A27").Value2
$r27_B = $ws.Range("B27").Value2
$r27_E = $ws.Range("E27").Value2
$r27_F = $ws.Range("F27").Value2
$r27_G = $ws.Range("G27").Value2
$r27_H = $ws.Range("H27").Value2
$r27_Q = $ws.Range("Q27").Value2
$r27_R = $ws.Range("R27").Value2
$r27_S = $ws.Range("S27").Value2
$r27_Z = $ws.Range("Z27").Value2
$r27_AB = $ws.Range("AB27").Value2

# --- row 25 <- old row 27 ---
$ws.Range("A25").Value = $r27_A
$ws.Range("B25").Value = $r27_B
$ws.Range("E25").Value = $r27_E
$ws.Range("F25").Value = $r27_F
$ws.Range("G25").Value = $r27_G
$ws.Range("H25").Value = $r27_H
$ws.Range("Q25").Value = $r27_Q
$ws.Range("R25").Value = $r27_R
$ws.Range("S25").Value = $r27_S
$ws.Range("Z25").Value = $r27_Z
$ws.Range("AB25").Value = $r27_AB

# --- row 26 <- old row 25 ---
$ws.Range("A26").Value = $r25_A
$ws.Range("B26").Value = $r25_B
$ws.Range("E26").Value = $r25_E
$ws.Range("F26").Value = $r25_F
$ws.Range("G26").Value = $r25_G
$ws.Range("H26").Value = $r25_H
$ws.Range("Q26").Value = $r25_Q
$ws.Range("R26").Value = $r25_R
$ws.Range("S26").Value = $r25_S
$ws.Range("Z26").Value = $r25_Z
$ws.Range("AB26").Value = $r25_AB

# --- row 27 <- old row 26 ---
$ws.Range("A27").Value = $r26_A
$ws.Range("B27").Value = $r26_B
$ws.Range("E27").Value = $r26_E
$ws.Range("F27").Value = $r26_F
$ws.Range("G27").Value = $r26_G
$ws.Range("H27").Value = $r26_H
$ws.Range("Q27").Value = $r26_Q
$ws.Range("R27").Value = $r26_R
$ws.Range("S27").Value = $r26_S
$ws.Range("Z27").Value = $r26_Z
$ws.Range("AB27").Value = $r26_AB
